$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''19.699.07'
$ws.Range("E2").Value = '  -8.99%  '
$ws.Range("D3").Value = '''1.384.34'
$ws.Range("E3").Value = '  -9.75%  '
$ws.Range("D4").Value = '''1.006'
$ws.Range("E4").Value = '  +0.52%  '
$ws.Range("E5").Value = '  +0.33%  '
$ws.Range("D6").Value = '''268.37'
$ws.Range("E6").Value = '  -6.92%  '
$ws.Range("D7").Value = '''0.3635'
$ws.Range("E7").Value = '  -7.87%  '
$ws.Range("D8").Value = '''0.3035'
$ws.Range("E8").Value = '  -4.06%  '
$ws.Range("D9").Value = '''38.15'
$ws.Range("E9").Value = '  -10.14%  '
$ws.Range("D10").Value = '''0.9699'
$ws.Range("E10").Value = '  -7.61%  '
$ws.Range("D11").Value = '''0.06389'
$ws.Range("E11").Value = '  -10.83%  '
$ws.Range("D12").Value = '''1.006'
$ws.Range("E12").Value = '  +0.53%  '
$ws.Range("D13").Value = '''5.270'
$ws.Range("E13").Value = '  -6.95%  '
$ws.Range("D14").Value = '''6.035'
$ws.Range("E14").Value = '  -8.36%  '
$ws.Range("B15").Value = 'Solana'
$ws.Range("C15").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D15").Value = '''16.34'
$ws.Range("E15").Value = '  -11.87%  '
$ws.Range("B16").Value = 'WrappedEther'
$ws.Range("C16").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D16").Value = '''1.385.90'
$ws.Range("E16").Value = '  -9.86%  '
$ws.Range("D17").Value = '''0.000009852'
$ws.Range("E17").Value = '  -9.53%  '
$ws.Range("D18").Value = '''0.05615'
$ws.Range("E18").Value = '  -14.91%  '
$ws.Range("E19").Value = '  +0.34%  '
$ws.Range("D21").Value = '''5.493'
$ws.Range("E21").Value = '  -10.13%  '
$ws.Range("D22").Value = '''14.26'
$ws.Range("E22").Value = '  -7.65%  '
$ws.Range("D23").Value = '''10.47'
$ws.Range("E23").Value = '  -2.67%  '
$ws.Range("D24").Value = '''2.240'
$ws.Range("E24").Value = '  -4.48%  '
$ws.Range("D25").Value = '''19.705.13'
$ws.Range("E25").Value = '  -8.95%  '
$ws.Range("D26").Value = '''2.151'
$ws.Range("E26").Value = '  -8.43%  '
$ws.Range("D27").Value = '''135.76'
$ws.Range("E27").Value = '  -8.96%  '
$ws.Range("D28").Value = '''16.49'
$ws.Range("E28").Value = '  -9.90%  '
$ws.Range("D29").Value = '''1.545.31'
$ws.Range("E29").Value = '  -11.61%  '
$ws.Range("D30").Value = '''107.27'
$ws.Range("E30").Value = '  -8.34%  '
$ws.Range("D31").Value = '''3.816'
$ws.Range("E31").Value = '  -21.20%  '
$ws.Range("D32").Value = '''5.193'
$ws.Range("E32").Value = '  -13.31%  '
$ws.Range("D33").Value = '''0.7818'
$ws.Range("E33").Value = '  -16.92%  '
$ws.Range("D34").Value = '''0.07568'
$ws.Range("E34").Value = '  -7.12%  '
$ws.Range("D35").Value = '''8.199'
$ws.Range("E35").Value = '  -3.54%  '
$ws.Range("E36").Value = '  +0.32%  '
$ws.Range("D37").Value = '''4.671'
$ws.Range("E37").Value = '  -9.26%  '
$ws.Range("D38").Value = '''0.05564'
$ws.Range("E38").Value = '  -7.33%  '
$ws.Range("D39").Value = '''0.02002'
$ws.Range("D40").Value = '''0.1864'
$ws.Range("E40").Value = '  -7.58%  '
$ws.Range("D41").Value = '''9.920'
$ws.Range("E41").Value = '  -9.20%  '
$ws.Range("D42").Value = '''1.288'
$ws.Range("E42").Value = '  -11.50%  '
$ws.Range("D43").Value = '''1.041'
$ws.Range("E43").Value = '  -11.59%  '
$ws.Range("D44").Value = '''3.457'
$ws.Range("E44").Value = '  -6.74%  '
$ws.Range("D45").Value = '''0.5136'
$ws.Range("E45").Value = '  -10.67%  '
$ws.Range("D46").Value = '''11.75'
$ws.Range("E46").Value = '  -9.34%  '
$ws.Range("D47").Value = '''0.4940'
$ws.Range("E47").Value = '  -9.98%  '
$ws.Range("D48").Value = '''107.62'
$ws.Range("E48").Value = '  -7.30%  '
$ws.Range("D49").Value = '''1.709'
$ws.Range("E49").Value = '  -8.84%  '
$ws.Range("D50").Value = '''1.005'
$ws.Range("E50").Value = '  +0.48%  '
$ws.Range("D51").Value = '''1.029'
$ws.Range("E51").Value = '  -11.69%  '
